$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1251.9166
$ws.Range("I28").Value = 534
$ws.Range("J28").Value = 3405.6667
$ws.Range("K28").Value = 534
$ws.Range("L28").Value = 3405.6667
$ws.Range("M28").Value = -49
$ws.Range("N28").Value = -4375.6667
$ws.Range("H48").Value = 1694.8823
$ws.Range("I48").Value = 985.3333
$ws.Range("J48").Value = 3397.8
$ws.Range("K48").Value = 2955.9999
$ws.Range("L48").Value = 10193.4
$ws.Range("M48").Value = -2663.9999
$ws.Range("N48").Value = -10777.4
$ws.Range("H56").Value = 1694.8823
$ws.Range("I56").Value = 985.3333
$ws.Range("J56").Value = 3397.8
$ws.Range("K56").Value = 2955.9999
$ws.Range("L56").Value = 10193.4
$ws.Range("M56").Value = -2421.9999
$ws.Range("N56").Value = -11261.4
$ws.Range("H62").Value = 7034.524
$ws.Range("I62").Value = 7336.3
$ws.Range("K62").Value = 7336.3
$ws.Range("M62").Value = -6712.3
$ws.Range("H65").Value = 7034.524
$ws.Range("I65").Value = 7336.3
$ws.Range("K65").Value = 36681.5
$ws.Range("M65").Value = -33561.5
$ws.Range("H92").Value = 250926
$ws.Range("J92").Value = 1637.5
$ws.Range("L92").Value = 1637.5
$ws.Range("N92").Value = -4133.5
$ws.Range("H97").Value = 10999
$ws.Range("J97").Value = 10999
$ws.Range("L97").Value = 32997
$ws.Range("N97").Value = -33989
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").ClearContents()
$ws.Range("H112").Value = 4677.364
$ws.Range("J112").Value = 4645.1
$ws.Range("L112").Value = 13935.3
$ws.Range("N112").Value = -16151.3
$ws.Range("H132").Value = 4009.2156
$ws.Range("J132").Value = 2938.3333
$ws.Range("L132").Value = 8814.999899999999
$ws.Range("N132").Value = -13874.9999
$ws.Range("H135").Value = 686.6923
$ws.Range("I135").Value = 358.77777
$ws.Range("K135").Value = 3228.99993
$ws.Range("M135").Value = -693.9999299999999
$ws.Range("H138").Value = 4826.6665
$ws.Range("I138").Value = 4826.6665
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 14479.9995
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -9339.999500000002
$ws.Range("N138").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 450
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 450
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 450
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -674
$ws.Range("H8").Value = 999
$ws.Range("I8").Value = 999
$ws.Range("K8").Value = 999
$ws.Range("M8").Value = -855
$ws.Range("H122").Value = 2090.9614
$ws.Range("I122").Value = 2044.6316
$ws.Range("J122").Value = 2216.7144
$ws.Range("K122").Value = 6133.8948
$ws.Range("L122").Value = 6650.1432
$ws.Range("M122").Value = -3683.8948
$ws.Range("N122").Value = -11550.1432

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 450
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 450
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -680
$ws.Range("H13").Value = 47499
$ws.Range("J13").Value = 47499
$ws.Range("L13").Value = 47499
$ws.Range("N13").Value = -47835
$ws.Range("H105").Value = 3201.5
$ws.Range("I105").Value = 3201.5
$ws.Range("K105").Value = 3201.5
$ws.Range("M105").Value = -1454.5
$ws.Range("H107").Value = 2440.3845
$ws.Range("I107").Value = 1561.5862
$ws.Range("K107").Value = 1561.5862
$ws.Range("M107").Value = 358.4138
$ws.Range("H134").Value = 2569.2666
$ws.Range("I134").Value = 1865.1052
$ws.Range("J134").Value = 6391.857
$ws.Range("K134").Value = 5595.3156
$ws.Range("L134").Value = 19175.571
$ws.Range("M134").Value = -3060.3156
$ws.Range("N134").Value = -24245.571

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1882.25
$ws.Range("J16").Value = 2669.5715
$ws.Range("L16").Value = 2669.5715
$ws.Range("N16").Value = -3243.5715
$ws.Range("H60").Value = 27109
$ws.Range("I60").Value = 4022.5
$ws.Range("K60").Value = 4022.5
$ws.Range("M60").Value = -3511.5
$ws.Range("H86").Value = 8499
$ws.Range("I86").Value = 8499
$ws.Range("K86").Value = 8499
$ws.Range("M86").Value = -7376
$ws.Range("H89").Value = 8499
$ws.Range("I89").Value = 8499
$ws.Range("K89").Value = 42495
$ws.Range("M89").Value = -36879
$ws.Range("H105").Value = 1074.2273
$ws.Range("I105").Value = 1074.2273
$ws.Range("K105").Value = 1074.2273
$ws.Range("M105").Value = 672.7727
$ws.Range("H107").Value = 2925.8518
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H113").Value = 1882.25
$ws.Range("J113").Value = 2669.5715
$ws.Range("L113").Value = 2669.5715
$ws.Range("N113").Value = -7009.5715
$ws.Range("H122").Value = 1954.8
$ws.Range("I122").Value = 2085.75
$ws.Range("K122").Value = 6257.25
$ws.Range("M122").Value = -3807.25
$ws.Range("H132").Value = 2106.6296
$ws.Range("I132").Value = 1599.8572
$ws.Range("K132").Value = 4799.571599999999
$ws.Range("M132").Value = -2269.571599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 188.28572
$ws.Range("I17").Value = 154.5
$ws.Range("J17").Value = 233.33333
$ws.Range("K17").Value = 463.5
$ws.Range("L17").Value = 699.99999
$ws.Range("M17").Value = -294.5
$ws.Range("N17").Value = -1037.99999
$ws.Range("H75").Value = 298.5
$ws.Range("J75").Value = 298
$ws.Range("L75").Value = 894
$ws.Range("N75").Value = -2890
$ws.Range("H78").Value = 298.5
$ws.Range("J78").Value = 298
$ws.Range("L78").Value = 2682
$ws.Range("N78").Value = -12666
$ws.Range("H108").Value = 3293.75
$ws.Range("I108").Value = 3293.75
$ws.Range("K108").Value = 9881.25
$ws.Range("M108").Value = -7001.25
$ws.Range("H129").Value = 2880.7144
$ws.Range("J129").Value = 3514.6667
$ws.Range("L129").Value = 10544.0001
$ws.Range("N129").Value = -20544.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7754
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 7754
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 7754
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -8294
$ws.Range("H73").Value = 7754
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 7754
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 7754
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -9626
$ws.Range("H102").Value = 3292
$ws.Range("I102").Value = 2198.0952
$ws.Range("J102").Value = 5589.2
$ws.Range("K102").Value = 2198.0952
$ws.Range("L102").Value = 5589.2
$ws.Range("M102").Value = -576.0952000000002
$ws.Range("N102").Value = -8833.200000000001
$ws.Range("H132").Value = 46136.523
$ws.Range("I132").Value = 50301.953
$ws.Range("K132").Value = 150905.859
$ws.Range("M132").Value = -148375.859

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3869.5715
$ws.Range("J100").Value = 4237.6
$ws.Range("L100").Value = 4237.6
$ws.Range("N100").Value = -5319.6
$ws.Range("H122").Value = 4216.5938
$ws.Range("J122").Value = 4821.0557
$ws.Range("L122").Value = 14463.1671
$ws.Range("N122").Value = -19363.1671
$ws.Range("H136").Value = 2527.8823
$ws.Range("I136").Value = 2029.625
$ws.Range("K136").Value = 6088.875
$ws.Range("M136").Value = -3538.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 84499.5
$ws.Range("I45").Value = 69999.5
$ws.Range("J45").Value = 91749.5
$ws.Range("K45").Value = 69999.5
$ws.Range("L45").Value = 91749.5
$ws.Range("M45").Value = -69508.5
$ws.Range("N45").Value = -92731.5
$ws.Range("H74").Value = 77662.664
$ws.Range("I74").Value = 19988
$ws.Range("J74").Value = 106500
$ws.Range("K74").Value = 19988
$ws.Range("L74").Value = 106500
$ws.Range("M74").Value = -19052
$ws.Range("N74").Value = -108372
$ws.Range("H77").Value = 77662.664
$ws.Range("I77").Value = 19988
$ws.Range("J77").Value = 106500
$ws.Range("K77").Value = 59964
$ws.Range("L77").Value = 319500
$ws.Range("M77").Value = -55284
$ws.Range("N77").Value = -328860
$ws.Range("H136").Value = 5549.409
$ws.Range("I136").Value = 5140.5884
$ws.Range("K136").Value = 15421.7652
$ws.Range("M136").Value = -12871.7652
